# Replace Estonian party abbreviations with full descriptive names.
# Applies to both standalone party-code cells and '+'-joined coalition strings,
# across every worksheet in the workbook.

$map = @{
    'EKD' = 'EKD - Coalition Party  (Eestil Koonderakond, EKD)'
    'ERL' = 'ERL - People''s Union  of Estonia (Eestimaa Rahvaliit, ERL)'
    'EÜR' = 'EÜRP - United People''s Party of Estonia (Eestimaa Ühendatud Rahvapartei, EÜRP)'
    'IL ' = 'IL - Pro Patria Union (Isamaaliit, IL)'
    'KeE' = 'KeE - Centre Party (Keskerakond, KeE)'
    'M -' = 'M - Moderates  (Mõõdukad, M)'
    'RE ' = 'RE - Reform Party   (, RE)'
    'RL ' = 'RL - People''s Union  (, RL)'
    'RP ' = 'RP - Res Publica (Res Publica, RP)'
    'ER ' = 'ER - Estonian Greens (Eestimaa Rohelised, ER)'
    'IRL' = 'IRL - Pro Patria and Res Publica Union (Isamaa ja Res Publica Liit, IRL)'
    'SDT' = 'SDTP - Social Democratic Labour Party (, SDTP)'
    'EKR' = 'EKRE - Conservative People’s Party (Eesti Konservatiivne Rahvaerakond, EKRE)'
    'EVA' = 'EVA - Estonian Free Party (Eesti Vabaerakond, EVA)'
    'SDE' = 'SDE - Social Democratic Party (, SDE)'
}


function Convert-PartyString($text) {
    if ($null -eq $text) { return $text }
    if (-not ($text -is [string])) { return $text }

    if ($map.ContainsKey($text)) {
        return $map[$text]
    }

    if ($text.Contains('+')) {
        $parts = $text.Split('+')
        $allKnown = $true
        foreach ($p in $parts) {
            if (-not $map.ContainsKey($p)) {
                $allKnown = $false
                break
            }
        }
        if ($allKnown) {
            $newParts = @()
            foreach ($p in $parts) {
                $newParts += $map[$p]
            }
            return [string]::Join('+', $newParts)
        }
    }

    return $null
}

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
            $val = $cell.Value2
            if ($val -is [string]) {
                $converted = Convert-PartyString $val
                if ($null -ne $converted) {
                    $cell.Value = $converted
                }
            }
        }
    }
}

